$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Swap the match-detail columns (F..V) between a handful of row pairs.
# Columns A..E (index / pais / torneio / temporada / data_partida) are left
# untouched on every row - only the match data itself (teams, scores, odds,
# timestamps and url) moves.
# ---------------------------------------------------------------------------
function Swap-MatchRows($rowA, $rowB) {
    for ($col = 6; $col -le 22; $col++) {
        $valA = $ws.Cells.Item($rowA, $col).Value2
        $valB = $ws.Cells.Item($rowB, $col).Value2
        $ws.Cells.Item($rowA, $col).Value2 = $valB
        $ws.Cells.Item($rowB, $col).Value2 = $valA
    }
}

Swap-MatchRows 4 5
Swap-MatchRows 26 27
Swap-MatchRows 28 29

# ---------------------------------------------------------------------------
# Append a new match row (row 52) after the last existing data row (51).
# Copy formatting from row 51 first so styles (bold/bordered index column,
# date-formatted data_partida column) carry over to the new row.
# ---------------------------------------------------------------------------
$ws.Range("A51:V51").Copy() | Out-Null
$ws.Range("A52:V52").PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(52, 1).Value2 = 51
$ws.Cells.Item(52, 2).Value2 = "italy"
$ws.Cells.Item(52, 3).Value2 = "serie-a"
$ws.Cells.Item(52, 4).Value2 = "2023-2024"
$ws.Cells.Item(52, 5).Value2 = 45195.86458333334
$ws.Cells.Item(52, 6).Value2 = "Juventus"
$ws.Cells.Item(52, 7).Value2 = 1
$ws.Cells.Item(52, 8).Value2 = "Lecce"
$ws.Cells.Item(52, 9).Value2 = 0
$ws.Cells.Item(52, 10).Value2 = 1.53
$ws.Cells.Item(52, 11).Value2 = "17/09/2023 01:02"
$ws.Cells.Item(52, 12).Value2 = 1.52
$ws.Cells.Item(52, 13).Value2 = "26/09/2023 20:36"
$ws.Cells.Item(52, 14).Value2 = 4.17
$ws.Cells.Item(52, 15).Value2 = "17/09/2023 01:02"
$ws.Cells.Item(52, 16).Value2 = 4.26
$ws.Cells.Item(52, 17).Value2 = "26/09/2023 20:44"
$ws.Cells.Item(52, 18).Value2 = 6.89
$ws.Cells.Item(52, 19).Value2 = "17/09/2023 01:02"
$ws.Cells.Item(52, 20).Value2 = 7.26
$ws.Cells.Item(52, 21).Value2 = "26/09/2023 20:25"
$ws.Cells.Item(52, 22).Value2 = "https://www.betexplorer.com/football/italy/serie-a/juventus-lecce/vqRq6h8F/"
